$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.806.36"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +3.86%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.662.84"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +1.86%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.997"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.28%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "566.98"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +6.11%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.44"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +2.70%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.997"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.21%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.612"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +6.02%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.662.95"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +1.65%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.85"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.51%  "

# Row 11
$ws.Range("E11").Value = "  +4.53%  "

# Row 12
$ws.Range("E12").Value = "  +7.02%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.343"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +3.10%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.109.11"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.99%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "60.716.14"
$ws.Range("D15").ClearFormats()

# Row 16
$ws.Range("E16").Value = "  +5.55%  "

# Row 17
$ws.Range("E17").Value = "  +4.33%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.647.86"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +1.40%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.54"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +2.85%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "342.76"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +2.45%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.43"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +2.96%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.36"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +2.44%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.05%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "66.36"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.49%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.441"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +4.95%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.19%  "

# Row 28
$ws.Range("E28").Value = "  +4.22%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0803"
$ws.Range("D29").ClearFormats()

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.997"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.13%  "

# Row 31
$ws.Range("E31").Value = "  +4.11%  "

# Row 32
$ws.Range("E32").Value = "  +4.42%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "159.73"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +3.99%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "19.23"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +1.60%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.906"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +8.88%  "

# Row 37
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.17"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +5.45%  "

# Row 38
$ws.Range("B38").Value = "Fetch.AI"
$ws.Range("C38").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.902"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +10.36%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "37.46"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +1.06%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "303.48"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +6.73%  "

# Row 42
$ws.Range("E42").Value = "  +1.89%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.04%  "

# Row 44
$ws.Range("E44").Value = "  +4.84%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.603"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +1.30%  "

# Row 46
$ws.Range("B46").Value = "Hedera"
$ws.Range("C46").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0545"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +2.86%  "

# Row 47
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "128.13"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +13.66%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "19.33"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +1.54%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "10.68"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.18%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0237"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +5.08%  "

# Row 51
$ws.Range("E51").Value = "  +4.44%  "
